$wb = $excel.ActiveWorkbook

# --- Rename sheet6 from "事業投資" to "債務" ---
$ws6 = $wb.Worksheets.Item(6)
$ws6.Name = "債務"

# --- Sheet5 ("保險") : rewrite as insurance data, column-major order ---
$ws5 = $wb.Worksheets.Item(5)

$ws5.Range("B1").Value = "company"
$ws5.Range("B2").Value = "富邦人壽"
$ws5.Range("B3").Value = "國泰人壽"
$ws5.Range("B4").Value = "富邦人壽"
$ws5.Range("B5").Value = "富邦人壽"

$ws5.Range("C1").Value = "name"
$ws5.Range("C2").Value = "安泰喬祿還本終身壽險"
$ws5.Range("C3").Value = "國泰人壽利率變動型年金甲型"
$ws5.Range("C4").Value = "增額終身壽險"
$ws5.Range("C5").Value = "增額終身壽險"

$ws5.Range("D1").Value = "owner"
$ws5.Range("D2").Value = "林郁方"
$ws5.Range("D3").Value = "林郁方"
$ws5.Range("D4").Value = "林郁方"
$ws5.Range("D5").Value = "黃昱華"

$ws5.Range("E1").Value = "property_category"
$ws5.Range("E2").Value = "insurance"
$ws5.Range("E3").Value = "insurance"
$ws5.Range("E4").Value = "insurance"
$ws5.Range("E5").Value = "insurance"

$ws5.Range("F1").Value = "category"
$ws5.Range("F2").Value = "normal"
$ws5.Range("F3").Value = "normal"
$ws5.Range("F4").Value = "normal"
$ws5.Range("F5").Value = "normal"

$ws5.Range("G1").Value = "date"
$ws5.Range("G2").Value = "2012-04-23"
$ws5.Range("G3").Value = "2012-04-23"
$ws5.Range("G4").Value = "2012-04-23"
$ws5.Range("G5").Value = "2012-04-23"

$ws5.Range("H1").Value = "legislator_name"
$ws5.Range("H2").Value = "林郁方"
$ws5.Range("H3").Value = "林郁方"
$ws5.Range("H4").Value = "林郁方"
$ws5.Range("H5").Value = "林郁方"

$ws5.Range("I1").Value = "legislator_id"
$ws5.Range("I2").Value = 716
$ws5.Range("I3").Value = 716
$ws5.Range("I4").Value = 716
$ws5.Range("I5").Value = 716

$ws5.Range("J1").Value = "source_file"
$ws5.Range("J2").Value = "tmp5c281"
$ws5.Range("J3").Value = "tmp5c281"
$ws5.Range("J4").Value = "tmp5c281"
$ws5.Range("J5").Value = "tmp5c281"

$ws5.Range("K1").Value = "index"
$ws5.Range("K2").Value = 102
$ws5.Range("K3").Value = 103
$ws5.Range("K4").Value = 104
$ws5.Range("K5").Value = 105

# column A (index) for sheet5 data rows
$ws5.Range("A2").Value = 102
$ws5.Range("A3").Value = 103
$ws5.Range("A4").Value = 104
$ws5.Range("A5").Value = 105

# --- Sheet6 ("債務") : rewrite as debt data, column-major order ---
$ws6.Range("B1").Value = "species"
$ws6.Range("B2").Value = "房屋貸款"
$ws6.Range("B3").Value = "親友借款"
$ws6.Range("B4").Value = "親友借款"

$ws6.Range("C1").Value = "debtor"
$ws6.Range("C2").Value = "林郁方"
$ws6.Range("C3").Value = "林郁方"
$ws6.Range("C4").Value = "林郁方"

$ws6.Range("D1").Value = "owner"
$ws6.Range("D2").Value = "國泰世華台北分行臺北市中正區博愛路"
$ws6.Range("D3").Value = "賴志威臺北市信義區松德路"
$ws6.Range("D4").Value = "張宏瑋臺北市内湖區大湖山莊街"

$ws6.Range("E1").Value = "total"
$ws6.Range("E2").Value = 6148905
$ws6.Range("E3").Value = 1600000
$ws6.Range("E4").Value = 1700000

$ws6.Range("F1").Value = "register_date"
$ws6.Range("F2").Value = "97年03月26日"
$ws6.Range("F3").Value = "97年05月05日"
$ws6.Range("F4").Value = "97年05月15日"

$ws6.Range("G1").Value = "register_reason"
$ws6.Range("G2").Value = "房屋貸款"
$ws6.Range("G3").Value = "購買房屋自備款"
$ws6.Range("G4").Value = "購買房屋自備款"

$ws6.Range("H1").Value = "property_category"
$ws6.Range("H2").Value = "debt"
$ws6.Range("H3").Value = "debt"
$ws6.Range("H4").Value = "debt"

$ws6.Range("I1").Value = "category"
$ws6.Range("I2").Value = "normal"
$ws6.Range("I3").Value = "normal"
$ws6.Range("I4").Value = "normal"

$ws6.Range("J1").Value = "date"
$ws6.Range("J2").Value = "2012-04-23"
$ws6.Range("J3").Value = "2012-04-23"
$ws6.Range("J4").Value = "2012-04-23"

$ws6.Range("K1").Value = "legislator_name"
$ws6.Range("K2").Value = "林郁方"
$ws6.Range("K3").Value = "林郁方"
$ws6.Range("K4").Value = "林郁方"

$ws6.Range("L1").Value = "legislator_id"
$ws6.Range("L2").Value = 716
$ws6.Range("L3").Value = 716
$ws6.Range("L4").Value = 716

$ws6.Range("M1").Value = "source_file"
$ws6.Range("M2").Value = "tmp5c281"
$ws6.Range("M3").Value = "tmp5c281"
$ws6.Range("M4").Value = "tmp5c281"

$ws6.Range("N1").Value = "index"
$ws6.Range("N2").Value = 115
$ws6.Range("N3").Value = 116
$ws6.Range("N4").Value = 117

# column A (index) for sheet6 data rows
$ws6.Range("A2").Value = 115
$ws6.Range("A3").Value = 116
$ws6.Range("A4").Value = 117
